$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1247
$ws1.Range("F4").Value = 17081
$ws1.Range("F5").Value = 47
$ws1.Range("F9").Value = 1024
$ws1.Range("F13").Value = 11816
$ws1.Range("F14").Value = 31
$ws1.Range("F15").Value = 38
$ws1.Range("F16").Value = 1497
$ws1.Range("F17").Value = 4693
$ws1.Range("F19").Value = 49
$ws1.Range("F21").Value = 81
$ws1.Range("F22").Value = 916

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1247
$ws4.Range("F4").Value = 17081
$ws4.Range("F5").Value = 47
$ws4.Range("F9").Value = 1024
$ws4.Range("F15").Value = 11816
$ws4.Range("F16").Value = 31
$ws4.Range("F17").Value = 38
$ws4.Range("F18").Value = 1497
$ws4.Range("F19").Value = 4693
$ws4.Range("F21").Value = 49
$ws4.Range("F23").Value = 81
$ws4.Range("F24").Value = 916
